# Created experiment order generation script
# Rewrites the 5 task-order sheets (one per position, sheetId stays 1..5)
# with freshly generated names and stimulus-file orderings, matching a
# re-run of the participant's order-generation script.

$wb = $excel.ActiveWorkbook

function Set-TaskOrderSheet($ws, $newName, $values) {
    $ws.Name = $newName

    $dims = $ws.UsedRange.Rows.Count
    $oldDataRows = $dims - 1
    $newDataRows = $values.Length

    if ($newDataRows -lt $oldDataRows) {
        # Remove the now-unused trailing rows entirely.
        $startRow = $newDataRows + 2
        $endRow = $oldDataRows + 1
        $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 2)).EntireRow.Delete()
    } elseif ($newDataRows -gt $oldDataRows) {
        # Grow the table, copying column-A's numeric style down for new rows.
        for ($r = $oldDataRows + 2; $r -le $newDataRows + 1; $r++) {
            $ws.Cells.Item(2, 1).Copy($ws.Cells.Item($r, 1))
        }
    }

    for ($i = 0; $i -lt $values.Length; $i++) {
        $r = $i + 2
        $ws.Cells.Item($r, 1).Value = $i
        $ws.Cells.Item($r, 2).Value = $values[$i]
    }
}

# Position 1 (sheetId 1): vSAT task order
Set-TaskOrderSheet $wb.Worksheets.Item(1) "vSAT_TO-16515889839356427" @(
    "vSAT_stims-1651588983903915.csv",
    "SAT_stims-16515889838889134.csv",
    "SAT_stims-16515889838678377.csv",
    "vSAT_stims-16515889839189117.csv"
)

# Position 2 (sheetId 2): TOL task order
Set-TaskOrderSheet $wb.Worksheets.Item(2) "TOL_TO-1651588983984406" @(
    "MM_stims-1651588983951267.csv",
    "ZM_stims-16515889839356427.csv",
    "MM_stims-16515889839668906.csv",
    "ZM_stims-1651588983951267.csv",
    "MM_stims-16515889839825156.csv",
    "ZM_stims-16515889839668906.csv"
)

# Position 3 (sheetId 3): NB task order
Set-TaskOrderSheet $wb.Worksheets.Item(3) "NB_TO-1651588987176944" @(
    "ZB-match_5-1651588985046842.csv",
    "TB-16515889864990346.csv",
    "TB-16515889862940319.csv",
    "ZB-match_6-16515889845245242.csv",
    "TB-16515889871649098.csv",
    "OB-16515889855775435.csv",
    "OB-16515889853055665.csv",
    "ZB-match_6-16515889842974277.csv",
    "OB-16515889858301802.csv"
)

# Position 4 (sheetId 4): RS task order
Set-TaskOrderSheet $wb.Worksheets.Item(4) "RS_TO-16515889871799092" @(
    "eyes closed",
    "eyes open"
)

# Position 5 (sheetId 5): GNG task order
Set-TaskOrderSheet $wb.Worksheets.Item(5) "GNG_TO-16515889872249105" @(
    "go_stims-16515889871819103.csv",
    "GNG_stims-16515889872079067.csv",
    "go_stims-16515889872099085.csv",
    "GNG_stims-16515889872239084.csv"
)
